$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Tan Ah Kow"
$ws.Range("C12").Value = 50
$ws.Range("D12").Value = $true
